$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "46.032.69"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.86%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.593.37"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.52"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.50"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.73%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.598"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.66%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.579"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.96"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.40"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0838"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("E13").Value = "  -1.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.993.66"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.55%  "
$ws.Range("E15").Value = "  +1.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.618.61"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.913"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.83"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "46.188.94"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.97%  "
$ws.Range("E20").Value = "  +0.67%  "
$ws.Range("E21").Value = "  -4.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.69"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "295.00"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +14.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.75"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.06"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.23"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "29.75"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +4.54%  "
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.05"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.71"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.20"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "37.90"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.21"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.57"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.74%  "
$ws.Range("E35").Value = "  +3.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0837"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.41%  "
$ws.Range("E37").Value = "  -4.23%  "
$ws.Range("E38").Value = "  -6.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.121"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.95%  "
$ws.Range("E40").Value = "  +0.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.38"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +17.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "15.76"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0329"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.57"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.96"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -5.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.090.49"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "98.31"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +7.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.999"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.69"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +5.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "108.47"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.93%  "
$ws.Range("E51").Value = "  -0.62%  "
